# Auto-generated edit script: applies updated market-board / profit values
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 21057.074
$ws.Range("I15").Value = 21057.074
$ws.Range("K15").Value = 63171.222
$ws.Range("M15").Value = -63002.222
$ws.Range("H17").Value = 1953.4348
$ws.Range("J17").Value = 1953.4348
$ws.Range("L17").Value = 5860.3044
$ws.Range("N17").Value = -6196.3044
$ws.Range("H40").Value = 4764698.5
$ws.Range("I40").Value = 1719.3334
$ws.Range("K40").Value = 1719.3334
$ws.Range("M40").Value = -1544.3334
$ws.Range("H87").Value = 55499.5
$ws.Range("J87").Value = 55499.5
$ws.Range("L87").Value = 55499.5
$ws.Range("N87").Value = -57995.5
$ws.Range("H90").Value = 55499.5
$ws.Range("J90").Value = 55499.5
$ws.Range("L90").Value = 166498.5
$ws.Range("N90").Value = -178978.5
$ws.Range("H112").Value = 5308.857
$ws.Range("I112").Value = 1100
$ws.Range("K112").Value = 3300
$ws.Range("M112").Value = -2192
$ws.Range("H138").Value = 3552.69
$ws.Range("J138").Value = 4414.098
$ws.Range("L138").Value = 13242.294
$ws.Range("N138").Value = -23522.294

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2642396
$ws.Range("I32").Value = 3130264.8
$ws.Range("K32").Value = 3130264.8
$ws.Range("M32").Value = -3129977.8
$ws.Range("H61").Value = 31253432
$ws.Range("I61").Value = 2638.4119
$ws.Range("J61").Value = 66671000
$ws.Range("K61").Value = 2638.4119
$ws.Range("L61").Value = 66671000
$ws.Range("M61").Value = -2426.4119
$ws.Range("N61").Value = -66671424
$ws.Range("H136").Value = 31253432
$ws.Range("I136").Value = 2638.4119
$ws.Range("J136").Value = 66671000
$ws.Range("K136").Value = 7915.2357
$ws.Range("L136").Value = 200013000
$ws.Range("M136").Value = -5365.2357
$ws.Range("N136").Value = -200018100

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 8369324.5
$ws.Range("I86").Value = 10460346
$ws.Range("J86").Value = 5239.1665
$ws.Range("K86").Value = 10460346
$ws.Range("L86").Value = 5239.1665
$ws.Range("M86").Value = -10459223
$ws.Range("N86").Value = -7485.1665
$ws.Range("H89").Value = 8369324.5
$ws.Range("I89").Value = 10460346
$ws.Range("J89").Value = 5239.1665
$ws.Range("K89").Value = 52301730
$ws.Range("L89").Value = 26195.8325
$ws.Range("M89").Value = -52296114
$ws.Range("N89").Value = -37427.8325
$ws.Range("H94").Value = 2442.0588
$ws.Range("I94").Value = 1139.0416
$ws.Range("K94").Value = 1139.0416
$ws.Range("M94").Value = -688.0416
$ws.Range("H99").Value = 16670082
$ws.Range("I99").Value = 3995
$ws.Range("J99").Value = 25003124
$ws.Range("K99").Value = 3995
$ws.Range("L99").Value = 25003124
$ws.Range("M99").Value = -2497
$ws.Range("N99").Value = -25006120
$ws.Range("H105").Value = 3810.8708
$ws.Range("I105").Value = 2885.7334
$ws.Range("K105").Value = 2885.7334
$ws.Range("M105").Value = -1138.7334
$ws.Range("H109").Value = 59343
$ws.Range("J109").Value = 59343
$ws.Range("L109").Value = 59343
$ws.Range("N109").Value = -62117
$ws.Range("H112").Value = 39237.332
$ws.Range("J112").Value = 39237.332
$ws.Range("L112").Value = 39237.332
$ws.Range("N112").Value = -42191.332
$ws.Range("H134").Value = 5322421
$ws.Range("J134").Value = 7205.1113
$ws.Range("L134").Value = 21615.3339
$ws.Range("N134").Value = -26685.3339

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 874.2
$ws.Range("I22").Value = 843.5
$ws.Range("J22").Value = 997
$ws.Range("K22").Value = 843.5
$ws.Range("L22").Value = 997
$ws.Range("M22").Value = -493.5
$ws.Range("N22").Value = -1697
$ws.Range("H31").Value = 6419.722
$ws.Range("I31").Value = 2250
$ws.Range("J31").Value = 7092.2583
$ws.Range("K31").Value = 2250
$ws.Range("L31").Value = 7092.2583
$ws.Range("M31").Value = -1955
$ws.Range("N31").Value = -7682.2583
$ws.Range("H34").Value = 6419.722
$ws.Range("I34").Value = 2250
$ws.Range("J34").Value = 7092.2583
$ws.Range("K34").Value = 2250
$ws.Range("L34").Value = 7092.2583
$ws.Range("M34").Value = -2048
$ws.Range("N34").Value = -7496.2583

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1197.3636
$ws.Range("I5").Value = 774.1111
$ws.Range("J5").Value = 1705.2667
$ws.Range("K5").Value = 2322.3333
$ws.Range("L5").Value = 5115.800099999999
$ws.Range("M5").Value = -2210.3333
$ws.Range("N5").Value = -5339.800099999999
$ws.Range("H12").Value = 1295.04
$ws.Range("I12").Value = 1889.7273
$ws.Range("J12").Value = 827.7857
$ws.Range("K12").Value = 5669.1819
$ws.Range("L12").Value = 2483.3571
$ws.Range("M12").Value = -5496.1819
$ws.Range("N12").Value = -2829.3571
$ws.Range("H55").Value = 7707984.5
$ws.Range("J55").Value = 10017499
$ws.Range("L55").Value = 30052497
$ws.Range("N55").Value = -30052851
$ws.Range("H68").Value = 33336652
$ws.Range("I68").Value = 66667530
$ws.Range("J68").Value = 22226360
$ws.Range("K68").Value = 200002590
$ws.Range("L68").Value = 66679080
$ws.Range("M68").Value = -200001779
$ws.Range("N68").Value = -66680702
$ws.Range("H70").Value = 994.1667
$ws.Range("I70").Value = 994.1667
$ws.Range("K70").Value = 2982.5001
$ws.Range("M70").Value = -2667.5001
$ws.Range("H71").Value = 33336652
$ws.Range("I71").Value = 66667530
$ws.Range("J71").Value = 22226360
$ws.Range("K71").Value = 600007770
$ws.Range("L71").Value = 200037240
$ws.Range("M71").Value = -600003714
$ws.Range("N71").Value = -200045352
$ws.Range("H73").Value = 994.1667
$ws.Range("I73").Value = 994.1667
$ws.Range("K73").Value = 2982.5001
$ws.Range("M73").Value = -1890.5001
$ws.Range("H75").Value = 37056790
$ws.Range("I75").Value = 66669804
$ws.Range("J75").Value = 25667174
$ws.Range("K75").Value = 200009412
$ws.Range("L75").Value = 77001522
$ws.Range("M75").Value = -200008414
$ws.Range("N75").Value = -77003518
$ws.Range("H78").Value = 37056790
$ws.Range("I78").Value = 66669804
$ws.Range("J78").Value = 25667174
$ws.Range("K78").Value = 600028236
$ws.Range("L78").Value = 231004566
$ws.Range("M78").Value = -600023244
$ws.Range("N78").Value = -231014550
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H92").Value = 25642054
$ws.Range("I92").Value = 1544.5
$ws.Range("J92").Value = 76923070
$ws.Range("K92").Value = 4633.5
$ws.Range("L92").Value = 230769210
$ws.Range("M92").Value = -3385.5
$ws.Range("N92").Value = -230771706
$ws.Range("H135").Value = 1197.3636
$ws.Range("I135").Value = 774.1111
$ws.Range("J135").Value = 1705.2667
$ws.Range("K135").Value = 6966.9999
$ws.Range("L135").Value = 15347.4003
$ws.Range("M135").Value = -4431.9999
$ws.Range("N135").Value = -20417.4003
$ws.Range("H137").Value = 99953.42999999999
$ws.Range("J137").Value = 96045
$ws.Range("L137").Value = 288135
$ws.Range("N137").Value = -298335

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 200021
$ws.Range("I38").Value = 200021
$ws.Range("K38").Value = 200021
$ws.Range("M38").Value = -199558
$ws.Range("H102").Value = 4706.923
$ws.Range("I102").Value = 4980.909
$ws.Range("J102").Value = 3200
$ws.Range("K102").Value = 4980.909
$ws.Range("L102").Value = 3200
$ws.Range("M102").Value = -3358.909
$ws.Range("N102").Value = -6444
$ws.Range("H113").Value = 4454.919
$ws.Range("I113").Value = 2520
$ws.Range("J113").Value = 6497.3335
$ws.Range("K113").Value = 2520
$ws.Range("L113").Value = 6497.3335
$ws.Range("M113").Value = -350
$ws.Range("N113").Value = -10837.3335
$ws.Range("H132").Value = 3697.9
$ws.Range("I132").Value = 2938.8
$ws.Range("J132").Value = 4457
$ws.Range("K132").Value = 8816.400000000001
$ws.Range("L132").Value = 13371
$ws.Range("M132").Value = -6286.400000000001
$ws.Range("N132").Value = -18431

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3128833.8
$ws.Range("I61").Value = 4547927
$ws.Range("K61").Value = 4547927
$ws.Range("M61").Value = -4547725
$ws.Range("H113").Value = 3128833.8
$ws.Range("I113").Value = 4547927
$ws.Range("K113").Value = 4547927
$ws.Range("M113").Value = -4545757
$ws.Range("H122").Value = 3641.682
$ws.Range("I122").Value = 2359.8333
$ws.Range("K122").Value = 7079.499899999999
$ws.Range("M122").Value = -4629.499899999999
$ws.Range("H136").Value = 7865.8477
$ws.Range("I136").Value = 2642.6365
$ws.Range("J136").Value = 14495.308
$ws.Range("K136").Value = 7927.9095
$ws.Range("L136").Value = 43485.924
$ws.Range("M136").Value = -5377.9095
$ws.Range("N136").Value = -48585.924
$ws.Range("H140").Value = 75988.164
$ws.Range("J140").Value = 75988.164
$ws.Range("L140").Value = 75988.164
$ws.Range("N140").Value = -86348.164

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2499
$ws.Range("I100").Value = 1998
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 3996
$ws.Range("L100").Value = 6000
$ws.Range("M100").Value = -3455
$ws.Range("N100").Value = -7082
$ws.Range("H122").Value = 111752.375
$ws.Range("I122").Value = 183981.31
$ws.Range("K122").Value = 551943.9299999999
$ws.Range("M122").Value = -549493.9299999999
$ws.Range("H126").Value = 1062.6923
$ws.Range("I126").Value = 1045.4286
$ws.Range("K126").Value = 3136.2858
$ws.Range("M126").Value = -666.2857999999997
$ws.Range("H132").Value = 5197.231
$ws.Range("I132").Value = 4994.6787
$ws.Range("K132").Value = 14984.0361
$ws.Range("M132").Value = -12454.0361
$ws.Range("H136").Value = 16560490
$ws.Range("I136").Value = 23257086
$ws.Range("J136").Value = 563067.4399999999
$ws.Range("K136").Value = 69771258
$ws.Range("L136").Value = 1689202.32
$ws.Range("M136").Value = -69768708
$ws.Range("N136").Value = -1694302.32

